$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove rows that were dropped from the export ---
$accountsToRemove = @("004551472", "005053939", "004643737", "004884046", "004444164")
foreach ($acct in $accountsToRemove) {
    $match = $ws.Columns.Item(1).Find($acct)
    if ($match -ne $null) {
        $match.EntireRow.Delete()
    }
}

# --- Insert new row for account 004204500 (EDWARD) right above account 004479287 (ANA) ---
$anaCell = $ws.Columns.Item(1).Find("004479287")
$newRow = $anaCell.Row
$ws.Rows.Item($newRow).Insert()
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("A" + $newRow).Value = "004204500"
$ws.Range("B" + $newRow).Value = "EDWARD"
$ws.Range("C" + $newRow).Value = 114542.9

# --- Update existing balances ---
$anaCell = $ws.Columns.Item(1).Find("004479287")
$ws.Range("C" + $anaCell.Row).Value = 100839.84

$roberioCell = $ws.Columns.Item(1).Find("004586209")
$ws.Range("C" + $roberioCell.Row).Value = 6500

# --- Insert new row for account 004119016 (HEMAT) right below account 004643746 (MARIO) ---
$marioCell = $ws.Columns.Item(1).Find("004643746")
$newRow2 = $marioCell.Row + 1
$ws.Rows.Item($newRow2).Insert()
$ws.Range("A" + $newRow2).NumberFormat = "@"
$ws.Range("A" + $newRow2).Value = "004119016"
$ws.Range("B" + $newRow2).Value = "HEMAT"
$ws.Range("C" + $newRow2).Value = 399.92
